# CFFN_QTR_FIN.xlsx update: add two new most-recent quarterly columns
# (2018-12-31 and 2018-09-30) to each of the three statements (Income
# Statement, Balance Sheet, Cash Flow) by inserting two new columns at D:E
# and shifting the existing quarters right. Also correct a handful of
# mis-keyed figures that shipped with the prior quarter's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank columns before column D; this shifts the existing
#        D:K data (and its cell formatting) right to F:M.
$ws.Columns("D:E").Insert()

# --- 2. Carry the number formatting of the (now-shifted) first data column
#        F into the two new columns D and E, row by row across the used
#        range, so the new cells render as dates / thousands exactly like
#        their neighbours. Done in three blocks (one per statement) so the
#        blank separator rows (36, 78) are not touched and stay absent from
#        sheetData, matching the original layout.
$ws.Range("F5:F35").Copy()
$ws.Range("D5:D35").PasteSpecial(-4122)
$ws.Range("F5:F35").Copy()
$ws.Range("E5:E35").PasteSpecial(-4122)

$ws.Range("F38:F77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("F38:F77").Copy()
$ws.Range("E38:E77").PasteSpecial(-4122)

$ws.Range("F80:F102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$ws.Range("F80:F102").Copy()
$ws.Range("E80:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 3. Give the two new columns a sensible width, matching their
#        right-hand neighbours (F, G) rather than Excel's narrow default.
$ws.Range("D1").EntireColumn.ColumnWidth = $ws.Range("F1").EntireColumn.ColumnWidth
$ws.Range("E1").EntireColumn.ColumnWidth = $ws.Range("G1").EntireColumn.ColumnWidth

# --- 4. Populate the new column D (period ending 2018-12-31 / most
#        recent quarter) and column E (period ending 2018-09-30) for every
#        row that carries data in the statements.
$dVals = @{
    7 = 43465
    8 = 82400
    9 = "NA"
    10 = "NA"
    12 = "NA"
    13 = 0
    14 = 0
    15 = 0
    17 = 30100
    18 = 52300
    20 = -21300
    21 = 34100
    22 = 0
    23 = 30900
    24 = 6600
    25 = 0
    26 = 24400
    27 = 24400
    28 = 0
    29 = "NA"
    30 = 0
    31 = 0
    32 = 21300
    33 = 24400
    34 = 0
    35 = 24400
    38 = 43465
    41 = 37700
    42 = 144500
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 96100
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 9303800
    57 = 68500
    58 = 0
    59 = 3400
    60 = 0
    61 = 0
    62 = 18500
    63 = 0
    64 = 0
    65 = 0
    66 = 7957900
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 174000
    73 = 0
    74 = 0
    75 = 0
    76 = 1345900
    77 = 0
    80 = 43465
    81 = 24400
    83 = 3100
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 14300
    91 = -2200
    92 = 0
    93 = 0
    94 = 81600
    96 = -65400
    97 = 0
    98 = 0
    99 = 0
    100 = -151100
    101 = 0
    102 = -55200
}

$eVals = @{
    7 = 43373
    8 = 77300
    9 = "NA"
    10 = "NA"
    12 = "NA"
    13 = 0
    14 = 0
    15 = 0
    17 = 27200
    18 = 50100
    20 = -20900
    21 = 31600
    22 = 0
    23 = 29100
    24 = 7800
    25 = 0
    26 = 21400
    27 = 21400
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 20900
    33 = 21400
    34 = 0
    35 = 21400
    38 = 43373
    41 = 16300
    42 = 222500
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 96000
    49 = 0
    50 = 0
    51 = 0
    52 = 0
    53 = 0
    54 = 9449500
    57 = 83000
    58 = 0
    59 = 0
    60 = 0
    61 = 0
    62 = 21300
    63 = 0
    64 = 0
    65 = 0
    66 = 8057900
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 214600
    73 = 0
    74 = 0
    75 = 0
    76 = 1391600
    77 = 0
    80 = 43373
    81 = 21400
    83 = 2400
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = 34100
    91 = -4700
    92 = 0
    93 = 0
    94 = 20600
    96 = -11400
    97 = 0
    98 = 0
    99 = 0
    100 = -97700
    101 = 0
    102 = -43000
}

foreach ($r in $dVals.Keys) {
    $ws.Range("D$r").Value = $dVals[$r]
}
foreach ($r in $eVals.Keys) {
    $ws.Range("E$r").Value = $eVals[$r]
}

# --- 5. A few rows had the 2017-12-31 quarter (now column H after the
#        insert, previously column F) re-keyed with corrected figures at
#        the same time; fix those four cells.
$ws.Range("H24").Value = 8400
$ws.Range("H26").Value = 24300
$ws.Range("H27").Value = 24300
$ws.Range("H29").Value = 7500

$wb.Save()
